$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 110
$ws.Range("B2").Value = "CN0004"
$ws.Range("C2").Value = "INBOUND"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "Confirmed"
$ws.Range("G2").Value = 45376.14305555556
$ws.Range("G2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("G2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I2").Value = "NB001`n"

# Add row 3
$ws.Range("A3").Value = 118
$ws.Range("B3").Value = "CN0010"
$ws.Range("C3").Value = "INBOUND"
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = "Completed"
$ws.Range("H3").Value = 45376.14305555556
$ws.Range("H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I3").Value = "NB001`n"

# Add row 4
$ws.Range("A4").Value = 111
$ws.Range("B4").Value = "CN0005"
$ws.Range("C4").Value = "INBOUND"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = "Scheduled"
$ws.Range("I4").Value = "NB001`n"
